$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.766.19"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.336.55"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'580.54"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "'175.69"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "3.332.47"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  +6.02%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").Value = "'46.88"
$ws.Range("E12").Value = "  +4.73%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Value = "'690.98"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "3.872.90"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "'8.45"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "67.796.96"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'0.119"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "3.341.33"
$ws.Range("D20").Value = "'17.54"
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "'5.42"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("D24").Value = "'16.95"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "'100.79"
$ws.Range("E25").Value = "  +3.06%  "
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("D28").Value = "'9.53"
$ws.Range("E28").Value = "  +5.74%  "
$ws.Range("D29").Value = "'33.03"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "'8.56"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").Value = "'7.06"
$ws.Range("E31").Value = "  +7.36%  "
$ws.Range("D32").Value = "'567.38"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'11.00"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "3.705.16"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("D39").Value = "'35.23"
$ws.Range("E39").Value = "  +12.49%  "
$ws.Range("E40").Value = "  +5.11%  "
$ws.Range("D41").Value = "'3.17"
$ws.Range("E41").Value = "  +6.85%  "
$ws.Range("D42").Value = "'2.63"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("D43").Value = "0.0₃0673"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").Value = "'0.336"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("D45").Value = "'3.30"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "'2.66"
$ws.Range("E47").Value = "  +5.81%  "
$ws.Range("D48").Value = "'0.129"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "'131.94"
$ws.Range("E51").Value = "  +2.98%  "
